# Alpha/DB.xlsx -- "Add menu, modify dishes and ingr database and add ID taker"
#
# Sheet "Ингредиенты" (Ingredients, sheet1):
#   - insert a new ingredient row (id 4, "Курина грудка" / Chicken breast, price 60)
#     right before the trailing "*" note row, pushing that row down.
#
# Sheet "Блюда" (Dishes, sheet2):
#   - rework the first dish block (rename dish, flip the ID-taker cell to a
#     number, tweak quantities, freeze the computed totals to plain values)
#     and append a totals row below it.
#   - add two brand-new dish blocks ("fhfgh" and "dfglkdf") each with their
#     own ingredient rows and live SUM() totals.
#   - add a trailing "*" marker used as the ID taker.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Ингредиенты"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Push the "*" row (old row 6) down to row 7, freeing up row 6.
$ws1.Rows.Item(6).Insert()

# New row 6: chicken breast ingredient.
$ws1.Range("A6").Value = 4
$ws1.Range("B6").Value = "Курина грудка"
$ws1.Range("C6").Value = 60

# The freshly-inserted row picked up the column default formatting; match
# the workbook's existing convention (no explicit per-cell style) by
# copying the style off an already-unstyled neighbour in the same columns.
$ws1.Range("A6").Style = $ws1.Range("A5").Style
$ws1.Range("B6").Style = $ws1.Range("B5").Style
$ws1.Range("C6").Style = $ws1.Range("C5").Style
$ws1.Range("A7").Style = $ws1.Range("A5").Style

# ---------------------------------------------------------------------
# Sheet 2: "Блюда"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Shift the whole sheet down by one row, freeing up row 1 as a spacer
# above the first dish block (mirrors the blank rows used between blocks
# further down the sheet).
$ws2.Rows.Item(1).Insert()

# --- Dish block 1 (row 2): rename + tweak values -----------------------
$ws2.Range("A2").Value = "Помидоры с капустой"
$ws2.Range("H2").Value = 0

$ws2.Range("C4").Value = 10
$ws2.Range("D4").Value = 12
$ws2.Range("E4").Value = 0.33
$ws2.Range("F4").Value = 0.396

$ws2.Range("C5").Value = 44
$ws2.Range("E5").Value = 2.42
$ws2.Range("F5").Value = 2.97

# Totals row for dish block 1.
$ws2.Range("E7").Value = "Сумма: "
$ws2.Range("E8").Value = 2.75
$ws2.Range("F8").Value = 3.366

# --- Dish block 2 (rows 10-18): new dish "fhfgh" ------------------------
$ws2.Range("A10").Value = "fhfgh"
$ws2.Range("B10").Value = 4
$ws2.Range("G10").Value = "ID:"
$ws2.Range("H10").Value = 1

$ws2.Range("A11").Value = "Продукты"
$ws2.Range("B11").Value = "Цена за кг"
$ws2.Range("C11").Value = "Кол-во (сад)"
$ws2.Range("D11").Value = "Кол-во (ясли)"
$ws2.Range("E11").Value = "СумСт (сад)"
$ws2.Range("F11").Value = "СумСт (ясли)"

for ($r = 12; $r -le 15; $r++) {
    $ws2.Range("A$r").Value = "Помидоры"
    $ws2.Range("B$r").Value = 33
    $ws2.Range("C$r").Value = 0
    $ws2.Range("D$r").Value = 0
    $ws2.Range("E$r").Formula = "=B$r*C$r/1000"
    $ws2.Range("F$r").Formula = "=B$r*D$r/1000"
}

$ws2.Range("E17").Value = "Сумма: "
$ws2.Range("E18").Formula = "=SUM(E12:E15)"
$ws2.Range("F18").Formula = "=SUM(F12:F15)"

# --- Dish block 3 (rows 20-27): new dish "dfglkdf" ----------------------
$ws2.Range("A20").Value = "dfglkdf"
$ws2.Range("B20").Value = 3
$ws2.Range("G20").Value = "ID:"
$ws2.Range("H20").Value = 2

$ws2.Range("A21").Value = "Продукты"
$ws2.Range("B21").Value = "Цена за кг"
$ws2.Range("C21").Value = "Кол-во (сад)"
$ws2.Range("D21").Value = "Кол-во (ясли)"
$ws2.Range("E21").Value = "СумСт (сад)"
$ws2.Range("F21").Value = "СумСт (ясли)"

for ($r = 22; $r -le 24; $r++) {
    $ws2.Range("A$r").Value = "Помидоры"
    $ws2.Range("B$r").Value = 33
    $ws2.Range("C$r").Value = 0
    $ws2.Range("D$r").Value = 0
    $ws2.Range("E$r").Formula = "=B$r*C$r/1000"
    $ws2.Range("F$r").Formula = "=B$r*D$r/1000"
}

$ws2.Range("E26").Value = "Сумма: "
$ws2.Range("E27").Formula = "=SUM(E22:E24)"
$ws2.Range("F27").Formula = "=SUM(F22:F24)"

# --- ID taker ------------------------------------------------------------
$ws2.Range("H29").Value = "*"

# Leave the selection where the original author's session ended up.
[void]$ws2.Range("H29").Select()
